$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 193.92857
$ws.Range("J9").Value = 327.6
$ws.Range("L9").Value = 327.6
$ws.Range("N9").Value = -665.6
$ws.Range("H17").Value = 3288.889
$ws.Range("J17").Value = 3288.889
$ws.Range("L17").Value = 9866.667000000001
$ws.Range("N17").Value = -10202.667
$ws.Range("H43").Value = 30429
$ws.Range("I43").Value = 39399.75
$ws.Range("J43").Value = 24448.5
$ws.Range("K43").Value = 39399.75
$ws.Range("L43").Value = 24448.5
$ws.Range("M43").Value = -39330.75
$ws.Range("N43").Value = -24586.5
$ws.Range("H116").Value = 22770.691
$ws.Range("I116").Value = 4264.875
$ws.Range("K116").Value = 4264.875
$ws.Range("M116").Value = -822.875
$ws.Range("H135").Value = 4816
$ws.Range("I135").Value = 3645
$ws.Range("J135").Value = 9500
$ws.Range("K135").Value = 32805
$ws.Range("L135").Value = 85500
$ws.Range("M135").Value = -30270
$ws.Range("N135").Value = -90570
$ws.Range("H138").Value = 2407.5273
$ws.Range("I138").Value = 1276.1428
$ws.Range("K138").Value = 3828.4284
$ws.Range("M138").Value = 1311.5716

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3056.23
$ws.Range("I32").Value = 2454.042
$ws.Range("J32").Value = 14497.8
$ws.Range("K32").Value = 2454.042
$ws.Range("L32").Value = 14497.8
$ws.Range("M32").Value = -2167.042
$ws.Range("N32").Value = -15071.8
$ws.Range("H45").Value = 1871.1111
$ws.Range("I45").Value = 1789.3636
$ws.Range("K45").Value = 1789.3636
$ws.Range("M45").Value = -1412.3636
$ws.Range("H53").Value = 4300
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H61").Value = 6256581.5
$ws.Range("I61").Value = 7492.4614
$ws.Range("K61").Value = 7492.4614
$ws.Range("M61").Value = -7280.4614
$ws.Range("H136").Value = 6256581.5
$ws.Range("I136").Value = 7492.4614
$ws.Range("K136").Value = 22477.3842
$ws.Range("M136").Value = -19927.3842

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 194980
$ws.Range("J132").Value = 194980
$ws.Range("L132").Value = 194980
$ws.Range("N132").Value = -205100

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1316.2222
$ws.Range("I16").Value = 1437.0769
$ws.Range("J16").Value = 1002
$ws.Range("K16").Value = 1437.0769
$ws.Range("L16").Value = 1002
$ws.Range("M16").Value = -1150.0769
$ws.Range("N16").Value = -1576
$ws.Range("H31").Value = 1295303.2
$ws.Range("I31").Value = 2061248.5
$ws.Range("J31").Value = 2770.875
$ws.Range("K31").Value = 2061248.5
$ws.Range("L31").Value = 2770.875
$ws.Range("M31").Value = -2060953.5
$ws.Range("N31").Value = -3360.875
$ws.Range("H34").Value = 1295303.2
$ws.Range("I34").Value = 2061248.5
$ws.Range("J34").Value = 2770.875
$ws.Range("K34").Value = 2061248.5
$ws.Range("L34").Value = 2770.875
$ws.Range("M34").Value = -2061046.5
$ws.Range("N34").Value = -3174.875
$ws.Range("H53").Value = 55000
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H96").Value = 8068
$ws.Range("J96").Value = 8068
$ws.Range("L96").Value = 8068
$ws.Range("N96").Value = -13560
$ws.Range("H113").Value = 1316.2222
$ws.Range("I113").Value = 1437.0769
$ws.Range("J113").Value = 1002
$ws.Range("K113").Value = 1437.0769
$ws.Range("L113").Value = 1002
$ws.Range("M113").Value = 732.9231
$ws.Range("N113").Value = -5342
$ws.Range("H117").Value = 30673.5
$ws.Range("J117").Value = 40700
$ws.Range("L117").Value = 40700
$ws.Range("N117").Value = -49878
$ws.Range("H132").Value = 4942.5
$ws.Range("I132").Value = 4931
$ws.Range("K132").Value = 14793
$ws.Range("M132").Value = -12263

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 5253.067
$ws.Range("I44").Value = 449.33334
$ws.Range("J44").Value = 8455.556
$ws.Range("K44").Value = 1348.00002
$ws.Range("L44").Value = 25366.668
$ws.Range("M44").Value = -950.0000199999999
$ws.Range("N44").Value = -26162.668
$ws.Range("H50").Value = 210323.67
$ws.Range("I50").Value = 1924.1111
$ws.Range("J50").Value = 835522.3
$ws.Range("K50").Value = 5772.3333
$ws.Range("L50").Value = 2506566.9
$ws.Range("M50").Value = -5291.3333
$ws.Range("N50").Value = -2507528.9
$ws.Range("H53").Value = 210323.67
$ws.Range("I53").Value = 1924.1111
$ws.Range("J53").Value = 835522.3
$ws.Range("K53").Value = 5772.3333
$ws.Range("L53").Value = 2506566.9
$ws.Range("M53").Value = -5291.3333
$ws.Range("N53").Value = -2507528.9
$ws.Range("H101").Value = 11123333
$ws.Range("J101").Value = 11123333
$ws.Range("L101").Value = 33369999
$ws.Range("N101").Value = -33374867
$ws.Range("H107").Value = 4435.909
$ws.Range("I107").Value = 559.2
$ws.Range("K107").Value = 1677.6
$ws.Range("M107").Value = 242.3999999999999
$ws.Range("H108").Value = 396
$ws.Range("I108").Value = 396
$ws.Range("K108").Value = 1188
$ws.Range("M108").Value = 1692
$ws.Range("H109").Value = 1713.2858
$ws.Range("I109").Value = 1713.2858
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 5139.857400000001
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -4099.857400000001
$ws.Range("N109").ClearContents()
$ws.Range("H110").Value = 10358.5
$ws.Range("I110").Value = 10358.5
$ws.Range("K110").Value = 31075.5
$ws.Range("M110").Value = -26985.5
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H114").Value = 5703.143
$ws.Range("I114").Value = 3193.889
$ws.Range("J114").Value = 10219.8
$ws.Range("K114").Value = 9581.667000000001
$ws.Range("L114").Value = 30659.4
$ws.Range("M114").Value = -6327.667000000001
$ws.Range("N114").Value = -37167.39999999999
$ws.Range("H117").Value = 1648.6666
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H119").Value = 4858
$ws.Range("J119").Value = 15000
$ws.Range("L119").Value = 45000
$ws.Range("N119").Value = -54676
$ws.Range("H120").Value = 14007.125
$ws.Range("I120").Value = 14007.125
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 42021.375
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -37183.375
$ws.Range("N120").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 4767.273
$ws.Range("I18").Value = 3430
$ws.Range("K18").Value = 3430
$ws.Range("M18").Value = -3137

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 6104.143
$ws.Range("J9").Value = 20350
$ws.Range("L9").Value = 20350
$ws.Range("N9").Value = -20798
$ws.Range("H61").Value = 10881.632
$ws.Range("I61").Value = 9816.647000000001
$ws.Range("K61").Value = 9816.647000000001
$ws.Range("M61").Value = -9614.647000000001
$ws.Range("H113").Value = 10881.632
$ws.Range("I113").Value = 9816.647000000001
$ws.Range("K113").Value = 9816.647000000001
$ws.Range("M113").Value = -7646.647000000001
$ws.Range("H133").Value = 59896
$ws.Range("J133").Value = 59896
$ws.Range("L133").Value = 59896
$ws.Range("N133").Value = -64956

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3244.111
$ws.Range("I100").Value = 3149.625
$ws.Range("K100").Value = 6299.25
$ws.Range("M100").Value = -5758.25
$ws.Range("H135").Value = 125000
$ws.Range("J135").Value = 125000
$ws.Range("L135").Value = 125000
$ws.Range("N135").Value = -135140
